# The commit adds a new data row (id 12 / "jameel khan") at the bottom of
# the table on the active sheet, extending the used range from A1:E12 to
# A1:E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "jameel khan"
$ws.Range("C13").Value = "jameela@gmail.com"
$ws.Range("D13").Value = 9654879897
$ws.Range("E13").Value = "Delhi"
